$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 70; existing rows 70-94 shift down to 73-97.
$ws.Rows("70:72").Insert()

# Fill in the 3 new rows (weekly update for Femacal de La Calera - Membrillo).
# Row 70: Especial
$ws.Range("A70").Value = 3
$ws.Range("B70").Value = "Femacal de La Calera"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 45001
$ws.Range("E70").Value = 5
$ws.Range("F70").Value = "Fruta"
$ws.Range("G70").Value = 100104
$ws.Range("H70").Value = "Frutos de pepita"
$ws.Range("I70").Value = 100104003
$ws.Range("J70").Value = "Membrillo"
$ws.Range("K70").Value = "Champion"
$ws.Range("L70").Value = "Especial"
$ws.Range("M70").Value = 68
$ws.Range("N70").Value = 18000
$ws.Range("O70").Value = 18000
$ws.Range("P70").Value = 18000
$ws.Range("Q70").Value = "$/caja 18 kilos empedrada"
$ws.Range("R70").Value = "Región de O'Higgins"
$ws.Range("S70").Value = 1000
$ws.Range("T70").Value = 18

# Row 71: Primera
$ws.Range("A71").Value = 3
$ws.Range("B71").Value = "Femacal de La Calera"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 45001
$ws.Range("E71").Value = 5
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100104
$ws.Range("H71").Value = "Frutos de pepita"
$ws.Range("I71").Value = 100104003
$ws.Range("J71").Value = "Membrillo"
$ws.Range("K71").Value = "Champion"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 67
$ws.Range("N71").Value = 16000
$ws.Range("O71").Value = 16000
$ws.Range("P71").Value = 16000
$ws.Range("Q71").Value = "$/caja 18 kilos empedrada"
$ws.Range("R71").Value = "Región de O'Higgins"
$ws.Range("S71").Value = 889
$ws.Range("T71").Value = 18

# Row 72: Segunda
$ws.Range("A72").Value = 3
$ws.Range("B72").Value = "Femacal de La Calera"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 45001
$ws.Range("E72").Value = 5
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100104
$ws.Range("H72").Value = "Frutos de pepita"
$ws.Range("I72").Value = 100104003
$ws.Range("J72").Value = "Membrillo"
$ws.Range("K72").Value = "Champion"
$ws.Range("L72").Value = "Segunda"
$ws.Range("M72").Value = 50
$ws.Range("N72").Value = 14000
$ws.Range("O72").Value = 14000
$ws.Range("P72").Value = 14000
$ws.Range("Q72").Value = "$/caja 18 kilos empedrada"
$ws.Range("R72").Value = "Región de O'Higgins"
$ws.Range("S72").Value = 778
$ws.Range("T72").Value = 18
